$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.350.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.86%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.547.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.78%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.53%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.98"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.542.47"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.75%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.611"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.19%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.78%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.42"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +8.82%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.584"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.20%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.34"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.35%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000276"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.35%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.129.78"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.99%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.31"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.52%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "610.00"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.52%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.555.18"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.88%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.462.58"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.88%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.72%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.28"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.76%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.877"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.95%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.17"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -17.31%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.64"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.51"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.01%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.71"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.38%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.60"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.79%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.78"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.38%  "

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "672.00"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +9.12%  "

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.21"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.66%  "

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.03"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.11%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.07"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.13%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.65%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.60"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.98%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.100"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.14%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.73"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.22%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0475"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.63%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "57.35"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.44%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.12%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.26%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.374.70"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.82%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.318"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.66%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0698"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.10%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.83%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "32.58"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.80%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.59"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.17%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.51%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.31"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.28%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.06%  "
